$d = $word.ActiveDocument

# Locate the paragraph whose text is "7" (the start of the block to remove)
# and the paragraph whose text is "10" (the end of the block; only its text
# run is removed, the bookmark + paragraph itself survives).
$sevenIndex = -1
$tenIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t -eq "7") { $sevenIndex = $i }
    if ($t -eq "10") { $tenIndex = $i }
}

# The paragraph immediately before "7" (the blank separator paragraph added
# right after "6") is removed along with "7", "8" and "9" in one delete that
# spans from the start of that blank paragraph to the end of the "9" paragraph
# (i.e. right up to, but not including, the "10" paragraph).
$blankIndex = $sevenIndex - 1
$startDelete = $d.Paragraphs.Item($blankIndex).Range.Start
$endDelete = $d.Paragraphs.Item($tenIndex - 1).Range.End
$d.Range($startDelete, $endDelete).Delete()

# The "10" paragraph is now at index $blankIndex (everything before it shifted
# up). Clear just its text, keeping the bookmark and the paragraph mark.
$last = $d.Paragraphs.Item($blankIndex)
$lastRange = $last.Range
$textRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$textRange.Text = ""
